# Update countries & provincias Spain
# Applies the 24-May-2020 03:35 -> 04:05 COVID data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Numeric refresh: Corea del Sur (row 49) ---
$ws.Range("B49").Value = 11190
$ws.Range("C49").Value = 25
$ws.Range("D49").Value = 10213
$ws.Range("E49").Value = 711

# --- Sudan overtakes Honduras in the ranking (rows 73/74 swap) ---
# Row 73 becomes Sudan's (new, higher) figures
$ws.Range("A73").Value = "Sudan"
$ws.Range("B73").Value = 3628
$ws.Range("C73").Value = 250
$ws.Range("D73").Value = 424
$ws.Range("E73").Value = 3058
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 9
$ws.Range("H73").Value = 146

# Row 74 becomes Honduras, carrying its previous (unchanged) figures
$ws.Range("A74").Value = "Honduras"
$ws.Range("B74").Value = 3477
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 439
$ws.Range("E74").Value = 2871
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 167

# --- Vietnam numeric refresh (row 145) ---
$ws.Range("B145").Value = 325
$ws.Range("C145").Value = 1
$ws.Range("E145").Value = 58

# --- Belice overtakes Nueva Caledonia (rows 199/200 swap) ---
$ws.Range("A199").Value = "Belice"
$ws.Range("D199").Value = 16
$ws.Range("H199").Value = 2

$ws.Range("A200").Value = "Nueva Caledonia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

# --- Reorder tied entries: Sahara Occidental / San Bartolome / Bonaire... (rows 214-216) ---
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

# --- Refresh "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 04:05"
